$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking price strings so values
# like "1.000" / "0.4650" keep their exact text (not coerced to a number).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "27.871.14"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.756.09"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").Value = "327.17"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D7").Value = "0.4650"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "41.92"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "0.07354"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").Value = "1.079"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "20.49"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "5.977"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "7.141"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "1.752.48"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "92.07"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "0.00001053"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "0.06411"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "5.749"
$ws.Range("D23").Value = "27.887.54"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").Value = "2.160"
$ws.Range("E25").Value = "  +4.41%  "
$ws.Range("D26").Value = "161.78"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "20.08"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "1.956.67"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "2.137"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").Value = "122.84"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").Value = "1.063"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").Value = "0.09287"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "3.657"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "5.532"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.02268"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "11.63"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "0.06078"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").Value = "0.2060"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").Value = "4.895"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "0.6160"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").Value = "1.178"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "7.748"
$ws.Range("D43").Value = "1.340"
$ws.Range("E43").Value = "  -3.17%  "
$ws.Range("B44").Value = "PancakeSwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D44").Value = "3.731"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "12.97"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").Value = "0.5768"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").Value = "122.85"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").Value = "1.922"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "0.06792"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").Value = "1.118"
$ws.Range("D51").Value = "72.01"
$ws.Range("E51").Value = "  +0.04%  "
